# Weekly update of the "Albahaca" (Terminal La Palmera de La Serena) sheet.
# Inserts three new price records (weekly refresh) into the existing table:
#   - a new row before the current row 13 (date 2021-08-26)
#   - a new row before the current row 22/23 boundary (date 2021-08-27)
#   - a new row appended at the end of the table (date 2021-08-24)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fill-Row([int]$r, [double]$fecha, [double]$volumen, [double]$pmin, [double]$pmax, [double]$pprom) {
    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = 100112052
    $ws.Cells.Item($r, 7).Value = "Albahaca"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 9).Value = "Primera"
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $pmin
    $ws.Cells.Item($r, 12).Value = $pmax
    $ws.Cells.Item($r, 13).Value = $pprom
    $ws.Cells.Item($r, 14).Value = "`$/paquete"
    $ws.Cells.Item($r, 15).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

# 1) Insert new row at row 13 (pushes existing rows 13..25 down to 14..26)
$ws.Rows.Item(13).Insert()
Fill-Row 13 44434 600 4500 5000 4750

# 2) Insert new row at row 23 (pushes existing rows 23..26 down to 24..27)
$ws.Rows.Item(23).Insert()
Fill-Row 23 44435 1500 4500 5000 4750

# 3) Append new row at the end (row 28)
Fill-Row 28 44432 900 4500 5000 4750

$ws.Range("A1").Select()
